# Insert a new data row at row 310 (pushing existing rows 310-372 down to
# 311-373) and populate it with the new reading. Excel automatically
# extends the used range / dimension to A1:R373 and shifts the rest of the
# data down, which reproduces the diff (row 373 ends up holding what used
# to be row 372's values, etc.).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(310).Insert()

$ws.Range("A310").Value = 9
$ws.Range("B310").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C310").Value = "Metropolitana"
$ws.Range("D310").Value = 44641
$ws.Range("E310").Value = 13
$ws.Range("F310").Value = 100112039
$ws.Range("G310").Value = "Ciboulette"
$ws.Range("H310").Value = "Sin especificar"
$ws.Range("I310").Value = "Primera"
$ws.Range("J310").Value = 160
$ws.Range("K310").Value = 1600
$ws.Range("L310").Value = 1800
$ws.Range("M310").Value = 1700
$ws.Range("N310").Value = "$/docena de atados"
$ws.Range("O310").Value = "Región Metropolitana"
$ws.Range("P310").Value = 567
$ws.Range("Q310").Value = 3
$ws.Range("R310").Value = "Hortaliza"
